$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.3772241992882562
$ws1.Range("C2").Value = 0.06486486486486487
$ws1.Range("D2").Value = 0.8571428571428571
$ws1.Range("E2").Value = 0.1206030150753769
$ws1.Range("F2").Value = 0.2489626556016598
$ws1.Range("G2").Value = 0.5831775700934579
$ws1.Range("H2").Value = 0.7589954521134296
$ws1.Range("I2").Value = 24
$ws1.Range("J2").Value = 346
$ws1.Range("K2").Value = 188
$ws1.Range("L2").Value = 4

# --- Sheet: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")
$ws2.Range("B2").Value = 0.9791666666666666
$ws2.Range("C2").Value = 0.352059925093633
$ws2.Range("D2").Value = 0.5179063360881543

$ws2.Range("B3").Value = 0.06486486486486487
$ws2.Range("C3").Value = 0.8571428571428571
$ws2.Range("D3").Value = 0.1206030150753769

$ws2.Range("B4").Value = 0.3772241992882562
$ws2.Range("C4").Value = 0.3772241992882562
$ws2.Range("D4").Value = 0.3772241992882562
$ws2.Range("E4").Value = 0.3772241992882562

$ws2.Range("B5").Value = 0.5220157657657658
$ws2.Range("C5").Value = 0.6046013911182451
$ws2.Range("D5").Value = 0.3192546755817656

$ws2.Range("B6").Value = 0.933614263729922
$ws2.Range("C6").Value = 0.3772241992882563
$ws2.Range("D6").Value = 0.4981118645786208

# --- Sheet: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 188
$ws3.Range("C2").Value = 346
$ws3.Range("B3").Value = 4
$ws3.Range("C3").Value = 24
